$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("PairCorrInput")
$ws2 = $wb.Worksheets.Item("Information")

# Update the NumDays value
$ws1.Range("E2").Value = 535

# Update the CC Averages period values
$ws1.Range("F2").Value = 90
$ws1.Range("F3").Value = 180
$ws1.Range("F4").Value = 365
$ws1.Range("F5").Value = 525

# Add two new formatted (blank, date-formatted) cells below the existing K11/K12
# pair by copying the existing number format down into F21:F22
$ws1.Range("K11").Copy()
$ws1.Range("F21:F22").PasteSpecial(-4122)

# Restore selections on each sheet
$ws1.Range("D11").Select()
$ws2.Range("F33").Select()

# Information becomes the active/visible sheet
$ws2.Activate()
